# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the Hades profit-tracking workbook
# (sheets ALC, ARM, BSM, CRP, CUL, LTW, WVR) per the authoritative diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 578.9459000000001
$ws.Range("J17").Value = 578.9459000000001
$ws.Range("L17").Value = 1736.8377
$ws.Range("N17").Value = -2072.8377
$ws.Range("H69").Value = 2892
$ws.Range("I69").Value = 1800
$ws.Range("J69").Value = 3013.3333
$ws.Range("K69").Value = 5400
$ws.Range("L69").Value = 9039.999899999999
$ws.Range("M69").Value = -4526
$ws.Range("N69").Value = -10787.9999
$ws.Range("H72").Value = 2892
$ws.Range("I72").Value = 1800
$ws.Range("J72").Value = 3013.3333
$ws.Range("K72").Value = 16200
$ws.Range("L72").Value = 27119.9997
$ws.Range("M72").Value = -11832
$ws.Range("N72").Value = -35855.9997
$ws.Range("H74").Value = 3955.0264
$ws.Range("I74").Value = 3965.3438
$ws.Range("K74").Value = 3965.3438
$ws.Range("M74").Value = -3029.3438
$ws.Range("H77").Value = 3955.0264
$ws.Range("I77").Value = 3965.3438
$ws.Range("K77").Value = 19826.719
$ws.Range("M77").Value = -15146.719
$ws.Range("H92").Value = 627.1111
$ws.Range("I92").Value = 428.41177
$ws.Range("J92").Value = 4005
$ws.Range("K92").Value = 428.41177
$ws.Range("L92").Value = 4005
$ws.Range("M92").Value = 819.5882300000001
$ws.Range("N92").Value = -6501
$ws.Range("H94").Value = 2549
$ws.Range("I94").Value = 2549
$ws.Range("K94").Value = 2549
$ws.Range("M94").Value = -2098
$ws.Range("H100").Value = 2074.875
$ws.Range("I100").Value = 950.6667
$ws.Range("J100").Value = 5447.5
$ws.Range("K100").Value = 950.6667
$ws.Range("L100").Value = 5447.5
$ws.Range("M100").Value = -409.6667
$ws.Range("N100").Value = -6529.5
$ws.Range("H103").Value = 801.13635
$ws.Range("I103").Value = 512
$ws.Range("J103").Value = 1420.7142
$ws.Range("K103").Value = 1536
$ws.Range("L103").Value = 4262.142599999999
$ws.Range("M103").Value = -950
$ws.Range("N103").Value = -5434.142599999999
$ws.Range("H106").Value = 5131873.5
$ws.Range("I106").Value = 3974.5
$ws.Range("J106").Value = 66666664
$ws.Range("K106").Value = 3974.5
$ws.Range("L106").Value = 66666664
$ws.Range("M106").Value = -3343.5
$ws.Range("N106").Value = -66667926
$ws.Range("H113").Value = 3350.2942
$ws.Range("I113").Value = 3325.625
$ws.Range("J113").Value = 3372.2222
$ws.Range("K113").Value = 3325.625
$ws.Range("L113").Value = 3372.2222
$ws.Range("M113").Value = -71.625
$ws.Range("N113").Value = -9880.2222
$ws.Range("H132").Value = 1022380.56
$ws.Range("I132").Value = 1390.3422
$ws.Range("J132").Value = 4902143.5
$ws.Range("K132").Value = 4171.0266
$ws.Range("L132").Value = 14706430.5
$ws.Range("M132").Value = -1641.0266
$ws.Range("N132").Value = -14711490.5
$ws.Range("H137").Value = 2859855
$ws.Range("I137").Value = 4002229
$ws.Range("J137").Value = 3920.1
$ws.Range("K137").Value = 12006687
$ws.Range("L137").Value = 11760.3
$ws.Range("M137").Value = -12004137
$ws.Range("N137").Value = -16860.3
$ws.Range("H138").Value = 2899652.2
$ws.Range("J138").Value = 4766155.5
$ws.Range("L138").Value = 14298466.5
$ws.Range("N138").Value = -14308746.5
$ws.Range("H141").Value = 702
$ws.Range("I141").Value = 702
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2106
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 3074
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1410.49
$ws.Range("I32").Value = 1035.0741
$ws.Range("J32").Value = 3010.9473
$ws.Range("K32").Value = 1035.0741
$ws.Range("L32").Value = 3010.9473
$ws.Range("M32").Value = -748.0741
$ws.Range("N32").Value = -3584.9473
$ws.Range("H61").Value = 24439792
$ws.Range("I61").Value = 27806266
$ws.Range("K61").Value = 27806266
$ws.Range("M61").Value = -27806054
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H136").Value = 24439792
$ws.Range("I136").Value = 27806266
$ws.Range("K136").Value = 83418798
$ws.Range("M136").Value = -83416248

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7803.7144
$ws.Range("I86").Value = 10785.192
$ws.Range("K86").Value = 10785.192
$ws.Range("M86").Value = -9662.191999999999
$ws.Range("H89").Value = 7803.7144
$ws.Range("I89").Value = 10785.192
$ws.Range("K89").Value = 53925.95999999999
$ws.Range("M89").Value = -48309.95999999999
$ws.Range("H94").Value = 779.7
$ws.Range("I94").Value = 819.5
$ws.Range("K94").Value = 819.5
$ws.Range("M94").Value = -368.5
$ws.Range("H99").Value = 1118.4849
$ws.Range("I99").Value = 1171.1765
$ws.Range("K99").Value = 1171.1765
$ws.Range("M99").Value = 326.8235
$ws.Range("H116").Value = 40000
$ws.Range("J116").Value = 40000
$ws.Range("L116").Value = 40000
$ws.Range("N116").Value = -49178
$ws.Range("H134").Value = 4792.326
$ws.Range("I134").Value = 4556.1577
$ws.Range("J134").Value = 5914.125
$ws.Range("K134").Value = 13668.4731
$ws.Range("L134").Value = 17742.375
$ws.Range("M134").Value = -11133.4731
$ws.Range("N134").Value = -22812.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3711.1875
$ws.Range("I31").Value = 2892.6667
$ws.Range("J31").Value = 4202.3
$ws.Range("K31").Value = 2892.6667
$ws.Range("L31").Value = 4202.3
$ws.Range("M31").Value = -2597.6667
$ws.Range("N31").Value = -4792.3
$ws.Range("H34").Value = 3711.1875
$ws.Range("I34").Value = 2892.6667
$ws.Range("J34").Value = 4202.3
$ws.Range("K34").Value = 2892.6667
$ws.Range("L34").Value = 4202.3
$ws.Range("M34").Value = -2690.6667
$ws.Range("N34").Value = -4606.3
$ws.Range("H132").Value = 36677.35
$ws.Range("I132").Value = 25762.17
$ws.Range("J132").Value = 64647.5
$ws.Range("K132").Value = 77286.50999999999
$ws.Range("L132").Value = 193942.5
$ws.Range("M132").Value = -74756.50999999999
$ws.Range("N132").Value = -199002.5
$ws.Range("H134").Value = 29622.977
$ws.Range("I134").Value = 2741.0356
$ws.Range("J134").Value = 87522.53999999999
$ws.Range("K134").Value = 8223.106800000001
$ws.Range("L134").Value = 262567.62
$ws.Range("M134").Value = -5688.106800000001
$ws.Range("N134").Value = -267637.62

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 15874001
$ws.Range("J131").Value = 1089.1296
$ws.Range("L131").Value = 3267.3888
$ws.Range("N131").Value = -13347.3888

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 720.73334
$ws.Range("I46").Value = 674.63635
$ws.Range("J46").Value = 847.5
$ws.Range("K46").Value = 674.63635
$ws.Range("L46").Value = 847.5
$ws.Range("M46").Value = -486.63635
$ws.Range("N46").Value = -1223.5
$ws.Range("H68").Value = 1704.1428
$ws.Range("I68").Value = 1585.8
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 1585.8
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -836.8
$ws.Range("N68").Value = -3498
$ws.Range("H71").Value = 1704.1428
$ws.Range("I71").Value = 1585.8
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 7929
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -4185
$ws.Range("N71").Value = -17488
$ws.Range("H93").Value = 200
$ws.Range("I93").Value = 200
$ws.Range("K93").Value = 200
$ws.Range("M93").Value = 1048

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 8626
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 8626
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 8626
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -10498
$ws.Range("H77").Value = 8626
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 8626
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 25878
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -35238
$ws.Range("H100").Value = 67612.336
$ws.Range("I100").Value = 55931.332
$ws.Range("J100").Value = 85133.836
$ws.Range("K100").Value = 111862.664
$ws.Range("L100").Value = 170267.672
$ws.Range("M100").Value = -111321.664
$ws.Range("N100").Value = -171349.672
